$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D and E in this sheet store values as text (inlineStr) rather than
# numbers/percentages, so we force Text number format before assigning the
# new string values to avoid Excel auto-converting them to numeric values.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "261.55"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "0.93%"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "27.20"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "0.94%"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "4.708"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "0.42%"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "2.90%"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "6.718"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "0.61%"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.8499"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "-1.02%"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.9165"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "-1.51%"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "0.87%"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.04637"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "-3.16%"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07089"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "0.16%"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.03144"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "-0.27%"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.09053"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "-0.85%"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.001531"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "-0.18%"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0006168"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "1.84%"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.006129"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "2.02%"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.466"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "0.08%"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "2.179"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "0.63%"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "0.95%"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.080"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "-0.98%"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "0.28%"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "-0.36%"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "0.09%"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.03918"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "1.88%"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "-0.31%"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.004136"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "5.11%"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "-0.69%"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.01350"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "-11.76%"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.00005173"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "1.51%"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000751"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "0.10%"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.03591"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "-34.16%"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.1667"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "26.16%"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.00002102"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "0.10%"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "0.10%"

